# Add Uncertainty to Model
# Update the container_cost value in B2 (Sheet1). The Total_Cost formula in
# B8 (SUM(B2:B7)) will recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 2.0255432572060301

$excel.Calculate()
